$d = $word.ActiveDocument

# Locate the paragraph that ends with "Events (event sourcing). Backends. Peers. DIDs."
# and insert the new block of paragraphs right after it.
$rng = $d.Content
$rng.Find.Execute("Events (event sourcing). Backends. Peers. DIDs.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rng.Collapse(0) | Out-Null

# Alternating blank / content paragraphs to insert, in document order.
$newParagraphs = @(
    "",
    "Semiotic Layer: (PersistenceType, PersistenceSubject, PersistenceMember, PersistenceValue);",
    "",
    "Type Kind: Domain Service Handler. Domain signatures (domain / range: Subject Kind / Object Kind). Domain graph mappings context handler: function P(S) : O.",
    "",
    "Subject Kind: domain persistence resource types (employee).",
    "",
    "Member Kind: persistence members resource types (employee/salary;ARS).",
    "",
    "Value Kind: range resource types (salary;ARS).",
    "",
    "Reify Persistence as Relationship (Values as Relation Resources). Align domain / range with domains / primitive types (Member Kind, salary;ARS)."
)

foreach ($text in $newParagraphs) {
    $rng.InsertParagraphAfter() | Out-Null
    $rng.Collapse(0) | Out-Null
    $rng.Move(1, 1) | Out-Null
    if ($text -ne "") {
        $rng.InsertAfter($text) | Out-Null
        $rng.Collapse(0) | Out-Null
    }
}
